$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column indices corresponding to B, C, D, E, F, I, L, N
$cols = @(2, 3, 4, 5, 6, 9, 12, 14)

# New values for rows 2..25 (data row index 0..23), one row per array, matching $cols order
$newValues = @(
    @(21.24103426297824, 10.16730127147443, 6.783034899723307, 9.242673668494666, 36.50089427488131, 28.57412395558024, 10.51400227631765, 18.8490871953992),
    @(20.73311303431587, 9.571374779614219, 6.813626447635501, 9.258936321943976, 36.25053531943211, 28.56997137617331, 10.49248674702804, 18.92323688013348),
    @(20.42113794443177, 9.188186371296116, 6.833358383625877, 9.269554832211197, 36.1080944637791, 28.57515571999618, 10.48153583486193, 18.97075048805329),
    @(20.2941812516534, 9.027830497371921, 6.841638131437128, 9.274041532107022, 36.0529296647111, 28.57920677570074, 10.47764380748843, 18.99061338602422),
    @(20.27311685253606, 9.000954446605759, 6.84302741041934, 9.274796193342343, 36.04394478096584, 28.5799963123915, 10.47703205817566, 18.99394189045221),
    @(20.41942477799911, 9.186040563501928, 6.833469079853289, 9.269614694835154, 36.10733877281238, 28.57520251507649, 10.48148103274055, 18.97101633678426),
    @(21.06604568684337, 9.965490839614231, 6.793386155724211, 9.248149828323218, 36.41225763017729, 28.57108347763704, 10.50611626462329, 18.87424282729008),
    @(22.32462432897433, 11.35190367399198, 6.7222970804816, 9.211065352732868, 37.09748599449426, 28.6246029503037, 10.57223251628999, 18.70015472599594),
    @(23.23267697735594, 12.29779322933993, 6.674629660981503, 9.18685053130913, 37.65079809015307, 28.70170920612008, 10.63146909418593, 18.58171520939555),
    @(23.64006818330397, 12.7298900396777, 6.653931447734021, 9.176488174595951, 37.91253805315009, 28.74501282236989, 10.66068086674392, 18.52986777789281),
    @(23.79335975593911, 12.88978341919006, 6.646235152498682, 9.172657788457808, 38.0130218233763, 28.76259308836847, 10.67206313834754, 18.51052507747642),
    @(23.76039163344783, 12.85551321178835, 6.647886388050013, 9.173478571884855, 37.9913211994109, 28.75875431951428, 10.669597600867, 18.51467796135169),
    @(23.65270006588486, 12.74311925640064, 6.65329543016428, 9.176171171967679, 37.92077792404056, 28.74643548077312, 10.66161090337857, 18.52827062168772),
    @(23.58660385160472, 12.67378922067842, 6.656627071353978, 9.17783264969602, 37.87774408802389, 28.73904372955536, 10.65676038318122, 18.53663435059106),
    @(23.20592454904577, 12.26903137171105, 6.676002167679415, 9.187540851740048, 37.63388840599705, 28.69904472196469, 10.62960516940929, 18.58514432422562),
    @(22.97081540977798, 12.0165367899396, 6.688140599012551, 9.193663575672497, 37.4868091482114, 28.6766136714958, 10.61352278486969, 18.61542295431744),
    @(22.83506048648382, 11.87903912708894, 6.69521510944418, 9.197246693790493, 37.4031623237573, 28.66448649774854, 10.60448613667635, 18.63302969374399),
    @(22.78901089657312, 11.83215378823494, 6.697626360256736, 9.198470446069216, 37.37500620465536, 28.66051348478783, 10.60146330695492, 18.63902392427746),
    @(22.99589885933139, 12.04182693240933, 6.686838838090857, 9.193005439541666, 37.50236827021763, 28.67892133226497, 10.61521272094718, 18.61217995599822),
    @(23.68435947295143, 12.7762332278824, 6.651702819253302, 9.175377751775082, 37.94146166194658, 28.75002175493066, 10.6639481384353, 18.524270245999),
    @(24.12853783968321, 13.23470802046167, 6.629565122831426, 9.164402554231705, 38.23637478608262, 28.8033788326156, 10.69766437610747, 18.46851056266728),
    @(23.89204955847351, 12.99199486383818, 6.641304903198068, 9.170210406697766, 38.07827296492252, 28.77427153328342, 10.67950060750448, 18.49811594326404),
    @(22.98456044952696, 12.03039947664987, 6.687427065286801, 9.193302786375705, 37.49533114797842, 28.67787564433229, 10.61444804844028, 18.61364549512675),
    @(21.98631678544643, 10.99264676656228, 6.74072620691786, 9.220563821023314, 36.90309456335991, 28.60350039645794, 10.55245811483962, 18.7455812981164)
)

for ($i = 0; $i -lt $newValues.Count; $i++) {
    $r = $i + 2
    $rowVals = $newValues[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $c = $cols[$j]
        $ws.Cells.Item($r, $c).Value = $rowVals[$j]
    }
}
